# Deploy app for the 2023 season
# Updates the Fixtures sheet: refreshes team/result data for rows 2-23,
# clears the old Semi-Final / Grand-Final rows (24-25), and updates the
# saved view state (scroll position + selection).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Fixtures")

$ws.Range("B2").Value = "Moreton Bay United"
$ws.Range("C2").Value = "Away"
$ws.Range("D2").Value = 26
$ws.Range("E2").Value = 2
$ws.Range("G2").Value = "Y"

$ws.Range("B3").Value = "Redlands United"
$ws.Range("C3").Value = "Away"
$ws.Range("D3").Value = 4
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = "Y"

$ws.Range("B4").Value = "Sunshine Coast Wanderers"
$ws.Range("C4").Value = "Away"
$ws.Range("D4").Value = 14
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = "Y"

$ws.Range("B5").Value = "Brisbane City"
$ws.Range("C5").Value = "Away"
$ws.Range("D5").Value = 18
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = "Y"

$ws.Range("B6").Value = "Brisbane Roar Youth"
$ws.Range("C6").Value = "Away"
$ws.Range("D6").Value = 26
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = "Y"

$ws.Range("B7").Value = "Peninsula Power"
$ws.Range("C7").Value = "Away"
$ws.Range("D7").Value = 9
$ws.Range("E7").Value = 4
$ws.Range("G7").Value = "Y"

$ws.Range("B8").Value = "Lions FC"
$ws.Range("C8").Value = "Away"
$ws.Range("D8").Value = 23
$ws.Range("E8").Value = 4
$ws.Range("G8").Value = "Y"

$ws.Range("B9").Value = "Gold Coast United"
$ws.Range("C9").Value = "Home"
$ws.Range("D9").Value = 30
$ws.Range("E9").Value = 4
$ws.Range("G9").Value = "N"

$ws.Range("B10").Value = "Rochedale Rovers"
$ws.Range("C10").Value = "Away"
$ws.Range("D10").Value = 5
$ws.Range("E10").Value = 5
$ws.Range("G10").Value = "N"

$ws.Range("B11").Value = "Gold Coast Knights"
$ws.Range("C11").Value = "Home"
$ws.Range("D11").Value = 14
$ws.Range("E11").Value = 5
$ws.Range("G11").Value = "N"

$ws.Range("B12").Value = "Eastern Suburbs"
$ws.Range("C12").Value = "Away"
$ws.Range("D12").Value = 19
$ws.Range("E12").Value = 5
$ws.Range("G12").Value = "N"

$ws.Range("B13").Value = "Moreton Bay United"
$ws.Range("C13").Value = "Home"
$ws.Range("D13").Value = 4
$ws.Range("E13").Value = 6
$ws.Range("G13").Value = "N"

$ws.Range("B14").Value = "Redlands United"
$ws.Range("C14").Value = "Home"
$ws.Range("D14").Value = 11
$ws.Range("E14").Value = 6
$ws.Range("G14").Value = "N"

$ws.Range("B15").Value = "Sunshine Coast Wanderers"
$ws.Range("C15").Value = "Home"
$ws.Range("D15").Value = 18
$ws.Range("E15").Value = 6
$ws.Range("G15").Value = "N"

$ws.Range("B16").Value = "Brisbane City"
$ws.Range("C16").Value = "Home"
$ws.Range("D16").Value = 25
$ws.Range("E16").Value = 6
$ws.Range("G16").Value = "N"

$ws.Range("B17").Value = "Brisbane Roar Youth"
$ws.Range("C17").Value = "Away"
$ws.Range("D17").Value = 8
$ws.Range("E17").Value = 7
$ws.Range("G17").Value = "N"

$ws.Range("B18").Value = "Peninsula Power"
$ws.Range("C18").Value = "Away"
$ws.Range("D18").Value = 14
$ws.Range("E18").Value = 7
$ws.Range("G18").Value = "N"

$ws.Range("B19").Value = "Lions FC"
$ws.Range("C19").Value = "Home"
$ws.Range("D19").Value = 23
$ws.Range("E19").Value = 7
$ws.Range("G19").Value = "N"

$ws.Range("B20").Value = "Gold Coast United"
$ws.Range("C20").Value = "Away"
$ws.Range("D20").Value = 30
$ws.Range("E20").Value = 7
$ws.Range("G20").Value = "N"

$ws.Range("B21").Value = "Rochedale Rovers"
$ws.Range("C21").Value = "Home"
$ws.Range("D21").Value = 13
$ws.Range("E21").Value = 8
$ws.Range("G21").Value = "N"

$ws.Range("B22").Value = "Gold Coast Knights"
$ws.Range("C22").Value = "Away"
$ws.Range("D22").Value = 20
$ws.Range("E22").Value = 8
$ws.Range("G22").Value = "N"

$ws.Range("B23").Value = "Eastern Suburbs"
$ws.Range("C23").Value = "Home"
$ws.Range("D23").Value = 26
$ws.Range("E23").Value = 8
$ws.Range("G23").Value = "N"

# The old Semi-Final (row 24) and Grand Final (row 25) fixtures are gone
# for the 2023 season deployment - clear their contents (formatting, like
# the Period column's number format, is left intact).
$ws.Range("A24:G25").ClearContents()

# Restore the saved view/selection state.
$ws.Range("B23").Select()
